# Restore/update the "From" value for rule R20 (row 10) in the Sample Project
# Main sheet: cell C10 changes from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
